$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 0.765625
$ws.Range("B6").Value = 0.734375
$ws.Range("B7").Value = 0.703125
$ws.Range("B8").Value = 0.765625
$ws.Range("B9").Value = 0.6875
$ws.Range("B10").Value = 0.765625
$ws.Range("B11").Value = 0.734375
$ws.Range("B13").Value = 0.75
$ws.Range("B14").Value = 0.671875
$ws.Range("B15").Value = 0.765625
$ws.Range("B16").Value = 0.703125
$ws.Range("B17").Value = 0.65625
$ws.Range("B18").Value = 0.65625
$ws.Range("B19").Value = 0.65625
$ws.Range("B20").Value = 0.609375
$ws.Range("B21").Value = 0.59375
$ws.Range("B22").Value = 0.59375
$ws.Range("B23").Value = 0.59375
$ws.Range("B24").Value = 0.609375
$ws.Range("B25").Value = 0.609375
$ws.Range("B29").Value = 0.609375
$ws.Range("B30").Value = 0.609375
$ws.Range("B67").Value = 0.625
$ws.Range("B68").Value = 0.625
$ws.Range("B69").Value = 0.625
$ws.Range("B70").Value = 0.625
$ws.Range("B71").Value = 0.625
$ws.Range("B72").Value = 0.625
$ws.Range("B73").Value = 0.625
$ws.Range("B74").Value = 0.625
$ws.Range("B75").Value = 0.625
$ws.Range("B76").Value = 0.625
$ws.Range("B77").Value = 0.625
$ws.Range("B78").Value = 0.625
$ws.Range("B79").Value = 0.625
$ws.Range("B80").Value = 0.625
$ws.Range("A102").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("A103").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("B103").Value = 0.53125
$ws.Range("A104").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("B104").Value = 0.5
$ws.Range("A105").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("B105").Value = 0.640625
$ws.Range("A106").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("A107").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("B107").Value = 0.53125
$ws.Range("A108").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("B108").Value = 0.578125
$ws.Range("A109").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("B109").Value = 0.515625
$ws.Range("A110").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("B110").Value = 0.609375
$ws.Range("A111").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("B111").Value = 0.5
$ws.Range("A112").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("A113").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("B113").Value = 0.609375
$ws.Range("A114").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("B114").Value = 0.546875
$ws.Range("A115").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("B115").Value = 0.53125
$ws.Range("A116").Value = "<__main__.DisplayOutputs object at 0x7f70ac15e700>"
$ws.Range("B116").Value = 0.6470588235294118
